$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "<Host appBase=""webapps"" autoDeploy=""true"" name=""localhost"" unpackWARs=""true"">`n<Valve className=""org.apache.catalina.valves.AccessLogValve"" directory=""logs"" pattern=""%h %l %u %t &quot;%r&quot; %s %b"" prefix=""localhost_access_log"" suffix="".txt""/>`n</Host>`n"

$ws.Range("B13").Value = $newText
